$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.929.86'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.899.53'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7962'
$ws.Range('E5').Value = '  -5.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.63'
$ws.Range('E6').Value = '  +1.21%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -4.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.59'
$ws.Range('E9').Value = '  -4.00%  '
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08123'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.676'
$ws.Range('E12').Value = '  +7.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7738'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').Value = '1.869.92'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.91'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.202'
$ws.Range('E16').Value = '  +5.45%  '
$ws.Range('D17').Value = '29.898.34'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.01'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.08'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007781'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.231'
$ws.Range('E21').Value = '  +17.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '2.135.91'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1673'
$ws.Range('E25').Value = '  -4.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.495'
$ws.Range('E26').Value = '  +2.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.29'
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.084'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('E30').Value = '  +3.51%  '
$ws.Range('E31').Value = '  +2.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.495'
$ws.Range('E32').Value = '  +4.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05654'
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.096'
$ws.Range('E34').Value = '  +0.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.296'
$ws.Range('E35').Value = '  +1.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7478'
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9990'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  -3.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01941'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.790'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').Value = '1.170.21'
$ws.Range('E41').Value = '  +15.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '75.04'
$ws.Range('E42').Value = '  +3.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4444'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.968'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8543'
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.38'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9997'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.14'
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.892'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.523'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.006'
$ws.Range('E51').Value = '  +10.46%  '
